# Add the "2022-Q1" sheet (new quarterly fund-holdings detail) right before
# the "总计" (grand-totals) sheet, and update "总计" with a new leading row
# summarizing the 2022-Q1 figures.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$totalSheet = $sheets.Item("总计")

# Insert the new sheet right before "总计" (Add's first positional arg is
# the "Before" sheet) so the tab order becomes ..., 2021-Q4, 2022-Q1, 总计.
$ws = $sheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# IMPORTANT: worksheet handles returned by this host track *position*, not
# stable identity. Now that a sheet has been spliced in right before it, the
# old $totalSheet handle resolves to the newly-added sheet instead of "总计".
# Re-fetch every handle we still need, by name, after any operation that can
# shift sheet positions.
$totalSheet = $sheets.Item("总计")

# "2021-Q1" has the same row count (12 data rows) we need here, so borrow its
# header/index-column direct formatting (thin border + bold, centered) via a
# formats-only paste - this reuses the existing style index (s="2") instead
# of synthesizing a new, slightly-different one (setting Font/Borders/
# Alignment properties directly on a fresh cell in this host does not merge
# into one style reliably).
$templateSheet = $sheets.Item("2021-Q1")
$templateSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A13").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)

# Helper: write a value as genuine text (so numeric-looking strings such as
# fund codes "001479" or figures "10.92" keep their literal form / leading
# zeros instead of silently becoming numbers), while leaving the cell's
# style index untouched (resets to the workbook "Normal" style - plain, no
# direct formatting - exactly like the sibling quarter sheets).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# code, name, fund size, stock position, position ratio, held mkt value (亿元), position rank
$data = @(
    @("519029", "华夏稳增混合", "10.92", "92.99", "3.84", "0.4193", 10),
    @("001479", "中邮风格轮动灵活配置混合", "9.45", "62.17", "2.72", "0.2570", 10),
    @("001305", "九泰天富改革新动力混合A", "3.74", "88.86", "5.56", "0.2079", 8),
    @("001782", "九泰久益灵活配置混合A", "2.33", "94.33", "5.62", "0.1309", 8),
    @("206013", "鹏华宏观灵活配置混合", "5.63", "32.95", "1.58", "0.0890", 8),
    @("001844", "九泰久益灵活配置混合C", "1.47", "94.33", "5.62", "0.0826", 8),
    @("006973", "太平睿盈混合A", "8.69", "29.33", "0.86", "0.0747", 10),
    @("009912", "九泰天富改革新动力混合C", "0.59", "88.86", "5.56", "0.0328", 8),
    @("350005", "天治中国制造2025灵活配置混合", "0.70", "58.26", "3.89", "0.0272", 2),
    @("350007", "天治趋势精选混合", "0.45", "40.93", "5.05", "0.0227", 1),
    @("007669", "太平睿盈混合C", "2.21", "29.33", "0.86", "0.0190", 10),
    @("350008", "天治新消费灵活配置混合", "0.09", "42.25", "4.70", "0.0042", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $i

    Set-TextValue $ws.Cells.Item($row, 2) $rec[0]
    Set-TextValue $ws.Cells.Item($row, 3) $rec[1]
    Set-TextValue $ws.Cells.Item($row, 4) $rec[2]
    Set-TextValue $ws.Cells.Item($row, 5) $rec[3]
    Set-TextValue $ws.Cells.Item($row, 6) $rec[4]
    Set-TextValue $ws.Cells.Item($row, 7) $rec[5]

    $ws.Cells.Item($row, 8).Value = $rec[6]
}

$ws.Range("A1").Select()

# --- Update the "总计" sheet: insert a new leading data row for 2022-Q1 ---
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 1).Style = $totalSheet.Cells.Item(3, 1).Style

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 1.37

# Renumber the index column (A) for the rows that shifted down one place.
for ($r = 3; $r -le 5; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$totalSheet.Range("A1").Select()
